$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the newly commented / completed order row values
$ws.Range("C2").Value = "Toy"
$ws.Range("B3").Value = "Test"
$ws.Range("D3").Value = "Santas Workshop - Essentials Edition"

# Move the active selection to B4, matching the saved cursor position
$ws.Range("B4").Select()
